# Insert a new data row at row 84 (pushing existing rows 84-128 down to 85-129)
# and populate it with the new "1a (cosecha)" Camote record for Ñuble /
# Terminal Hortofrutícola Agro Chillán - Zapallo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(84).Insert()

$ws.Range("A84").Value = 7
$ws.Range("B84").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C84").Value = 'Ñuble'
$ws.Range("D84").Value = 44606
$ws.Range("E84").Value = 16
$ws.Range("F84").Value = 100112045
$ws.Range("G84").Value = 'Zapallo'
$ws.Range("H84").Value = 'Camote'
$ws.Range("I84").Value = '1a (cosecha)'
$ws.Range("J84").Value = 200
$ws.Range("K84").Value = 350
$ws.Range("L84").Value = 400
$ws.Range("M84").Value = 375
$ws.Range("N84").Value = '$/kilo (volumen en unidades)'
$ws.Range("O84").Value = 'Región de O''Higgins'
$ws.Range("P84").Value = 375
$ws.Range("Q84").Value = 1
$ws.Range("R84").Value = 'Hortaliza'
